# Protocolo CMD - atualização de minutas de protocolo
# Troca o nome do representante da Primeira Outorgante (AMA) no preâmbulo
# do protocolo: "João Paulo Salazar Dias" -> "Ana Sofia Rodrigues dos Reis Mota".

$d = $word.ActiveDocument

$d.Content.Find.Execute(
    "João Paulo Salazar Dias",
    $true,
    $false,
    $false,
    $false,
    $false,
    $true,
    1,
    $false,
    "Ana Sofia Rodrigues dos Reis Mota",
    2
) | Out-Null
